$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.056.56"
$ws.Range("E2").Value = "  +1.62%  "

$ws.Range("D3").Value = "2.418.29"
$ws.Range("E3").Value = "  +1.90%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.32"
$ws.Range("E5").Value = "  +1.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.63"
$ws.Range("E6").Value = "  +4.24%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("E8").Value = "  +1.20%  "

$ws.Range("D9").Value = "2.415.30"
$ws.Range("E9").Value = "  +1.69%  "

$ws.Range("E10").Value = "  +2.45%  "

$ws.Range("E11").Value = "  -0.70%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.38"
$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.05"
$ws.Range("E14").Value = "  +4.31%  "

$ws.Range("E15").Value = "  +7.20%  "

$ws.Range("D16").Value = "2.842.51"
$ws.Range("E16").Value = "  +2.17%  "

$ws.Range("D17").Value = "62.044.79"
$ws.Range("E17").Value = "  +1.71%  "

$ws.Range("D18").Value = "2.417.32"
$ws.Range("E18").Value = "  +1.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.10"
$ws.Range("E19").Value = "  +2.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.19"
$ws.Range("E20").Value = "  +1.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.60"
$ws.Range("E21").Value = "  -0.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.71"
$ws.Range("E22").Value = "  +0.84%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.33"
$ws.Range("E24").Value = "  +1.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.74"
$ws.Range("E25").Value = "  +5.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.95"
$ws.Range("E26").Value = "  +7.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "574.93"
$ws.Range("E27").Value = "  +15.39%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.537.06"
$ws.Range("E28").Value = "  +2.61%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("D30").Value = "0.0₃0940"
$ws.Range("E30").Value = "  +6.49%  "

$ws.Range("E31").Value = "  +5.72%  "

$ws.Range("E32").Value = "  +2.15%  "

$ws.Range("E33").Value = "  -0.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("E34").Value = "  +2.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("E35").Value = "  +3.98%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("E37").Value = "  +5.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.82"
$ws.Range("E38").Value = "  +2.34%  "

$ws.Range("E39").Value = "  +1.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.72"
$ws.Range("E40").Value = "  +0.87%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "150.18"
$ws.Range("E41").Value = "  +4.06%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.85"
$ws.Range("E42").Value = "  -2.66%  "

$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.35"
$ws.Range("E44").Value = "  +14.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "150.00"
$ws.Range("E45").Value = "  +3.74%  "

$ws.Range("E46").Value = "  +1.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0544"
$ws.Range("E47").Value = "  +4.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.28"
$ws.Range("E48").Value = "  +5.36%  "

$ws.Range("E49").Value = "  +2.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0923"
$ws.Range("E50").Value = "  +1.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0229"
$ws.Range("E51").Value = "  +2.65%  "
